{"js": "// Split the single run \"{m\" (in the \"{m:v.name}\" Heading-1 paragraph) into\n// two runs \"{\" and \"m\", and split the single run \"{m:\" (in the\n// \"{m:endfor}\" paragraph) into two runs \"{\" and \"m:\".\n// This mirrors the upstream parser switching to\n// TokenIteratorFieldRewriterSplit, which tokenizes the M2Doc \"{m:...}\"\n// field markers run-by-run instead of keeping \"{m\" / \"{m:\" as one run.\n\nfunction wrapPkg(innerParagraphXml) {\n  return (\n    '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + innerParagraphXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// --- Edit 1: \"{m:v.name}\" heading paragraph -> split leading \"{m\" run ---\nconst headingPara = paragraphs.items.filter((p) => false); // placeholder, replaced below\nlet target1 = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"{m:v.name}\") {\n    target1 = paragraphs.items[i];\n    break;\n  }\n}\nif (target1) {\n  const ooxml1 =\n    '<w:p w:rsidR=\"00052FB8\" w:rsidRDefault=\"006F5523\" w:rsidP=\"00727C85\">' +\n    '<w:pPr><w:pStyle w:val=\"Titre1\"/></w:pPr>' +\n    '<w:r w:rsidR=\"008D51EC\"><w:t>{</w:t></w:r>' +\n    '<w:r w:rsidR=\"008D51EC\"><w:t>m</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">:v.name}</w:t></w:r>' +\n    '</w:p>';\n  target1.getRange().insertOoxml(wrapPkg(ooxml1), \"Replace\");\n  await context.sync();\n}\n\n// --- Edit 2: \"{m:endfor}\" paragraph -> split leading \"{m:\" run ---\nlet target2 = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"{m:endfor}\") {\n    target2 = paragraphs.items[i];\n    break;\n  }\n}\nif (target2) {\n  const ooxml2 =\n    '<w:p w:rsidR=\"00C52979\" w:rsidRDefault=\"006F5523\" w:rsidP=\"00F5495F\">' +\n    '<w:r w:rsidR=\"001B2B2B\"><w:t>{</w:t></w:r>' +\n    '<w:r w:rsidR=\"001B2B2B\"><w:t>m:</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">endfor}</w:t></w:r>' +\n    '</w:p>';\n  target2.getRange().insertOoxml(wrapPkg(ooxml2), \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Split the single run \"{m\" (in the \"{m:v.name}\" Heading-1 paragraph) into\n# two runs \"{\" and \"m\", and split the single run \"{m:\" (in the\n# \"{m:endfor}\" paragraph) into two runs \"{\" and \"m:\".\n# This mirrors the upstream parser switching to\n# TokenIteratorFieldRewriterSplit, which tokenizes the M2Doc \"{m:...}\"\n# field markers run-by-run instead of keeping \"{m\" / \"{m:\" as one run.\n\n$d = $word.ActiveDocument\n\n$wNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n\n    # --- Edit 1: \"{m:v.name}\" heading paragraph -> split leading \"{m\" run ---\n    if ($t.StartsWith(\"{m:v.name}\") -and $p.Style.NameLocal -eq \"Heading 1\") {\n        $xml = '<w:p ' + $wNs + ' w:rsidR=\"00052FB8\" w:rsidRDefault=\"006F5523\" w:rsidP=\"00727C85\">' +\n               '<w:pPr><w:pStyle w:val=\"Titre1\"/></w:pPr>' +\n               '<w:r w:rsidR=\"008D51EC\"><w:t>{</w:t></w:r>' +\n               '<w:r w:rsidR=\"008D51EC\"><w:t>m</w:t></w:r>' +\n               '<w:r><w:t xml:space=\"preserve\">:v.name}</w:t></w:r>' +\n               '</w:p>'\n        $p.Range.InsertXML($xml)\n    }\n\n    # --- Edit 2: \"{m:endfor}\" paragraph -> split leading \"{m:\" run ---\n    if ($t.StartsWith(\"{m:endfor}\")) {\n        $xml = '<w:p ' + $wNs + ' w:rsidR=\"00C52979\" w:rsidRDefault=\"006F5523\" w:rsidP=\"00F5495F\">' +\n               '<w:r w:rsidR=\"001B2B2B\"><w:t>{</w:t></w:r>' +\n               '<w:r w:rsidR=\"001B2B2B\"><w:t>m:</w:t></w:r>' +\n               '<w:r><w:t xml:space=\"preserve\">endfor}</w:t></w:r>' +\n               '</w:p>'\n        $p.Range.InsertXML($xml)\n    }\n}\n"}
